$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.354.27"
$ws.Range("E2").Value = "  +0.80%  "

$ws.Range("D3").Value = "2.411.26"
$ws.Range("E3").Value = "  -0.34%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "560.69"
$ws.Range("E5").Value = "  +1.80%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "135.87"
$ws.Range("E6").Value = "  -1.15%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.07%  "

$ws.Range("E8").Value = "  +0.54%  "

$ws.Range("E9").Value = "  +0.73%  "

$ws.Range("E10").Value = "  -0.79%  "

$ws.Range("E11").Value = "  +1.44%  "

$ws.Range("E12").Value = "  -1.44%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "24.71"
$ws.Range("E13").Value = "  -2.97%  "

$ws.Range("D14").Value = "2.839.07"
$ws.Range("E14").Value = "  -0.47%  "

$ws.Range("D15").Value = "60.258.23"
$ws.Range("E15").Value = "  +0.76%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000139"
$ws.Range("E16").Value = "  +0.87%  "

$ws.Range("D17").Value = "2.401.17"
$ws.Range("E17").Value = "  -1.51%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.23"
$ws.Range("E18").Value = "  -1.22%  "

$ws.Range("E19").Value = "  +3.40%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "324.75"
$ws.Range("E20").Value = "  -1.67%  "

$ws.Range("E21").Value = "  +1.80%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.00"
$ws.Range("E22").Value = "  -0.03%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "64.57"
$ws.Range("E23").Value = "  -3.11%  "

$ws.Range("E24").Value = "  +1.52%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.52"
$ws.Range("E25").Value = "  -2.52%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  -0.19%  "

$ws.Range("E27").Value = "  +0.67%  "

$ws.Range("E28").Value = "  +2.40%  "

$ws.Range("E29").Value = "  -0.35%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "170.65"
$ws.Range("E30").Value = "  +0.84%  "

$ws.Range("E31").Value = "  +0.20%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.08"
$ws.Range("E32").Value = "  +6.23%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.403"
$ws.Range("E33").Value = "  -1.94%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "18.36"
$ws.Range("E34").Value = "  -1.81%  "

$ws.Range("E35").Value = "  +3.70%  "

$ws.Range("E36").Value = "  +0.03%  "

$ws.Range("E37").Value = "  +0.04%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.19"
$ws.Range("E38").Value = "  -0.54%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "324.42"
$ws.Range("E39").Value = "  +3.27%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.61"
$ws.Range("E40").Value = "  -0.13%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "38.74"
$ws.Range("E41").Value = "  -2.01%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "147.73"
$ws.Range("E42").Value = "  +5.99%  "

$ws.Range("E43").Value = "  -2.90%  "

$ws.Range("E44").Value = "  +0.49%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "19.90"
$ws.Range("E45").Value = "  +1.71%  "

$ws.Range("E46").Value = "  -0.71%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.575"
$ws.Range("E47").Value = "  -0.49%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0222"
$ws.Range("E48").Value = "  -1.26%  "

$ws.Range("E49").Value = "  -0.01%  "

$ws.Range("E50").Value = "  +0.02%  "

$ws.Range("E51").Value = "  -0.59%  "
